$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.157.01'
$ws.Range("E2").Value = '  -0.64%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.587.88'
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.35'
$ws.Range("E5").Value = '  +0.56%  '

$ws.Range("E6").Value = '  -0.91%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  -0.18%  '

$ws.Range("E9").Value = '  -1.63%  '

$ws.Range("E10").Value = '  -2.45%  '

$ws.Range("E11").Value = '  -0.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.811.98'
$ws.Range("E12").Value = '  +0.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.581.94'
$ws.Range("E13").Value = '  -0.19%  '

$ws.Range("E14").Value = '  -1.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.509'
$ws.Range("E15").Value = '  -1.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.43'
$ws.Range("E16").Value = '  -1.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.183.51'
$ws.Range("E17").Value = '  -0.55%  '

$ws.Range("E18").Value = '  -0.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.36'
$ws.Range("E19").Value = '  -1.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '213.19'
$ws.Range("E20").Value = '  +1.04%  '

$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.22'
$ws.Range("E22").Value = '  -0.95%  '

$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("E24").Value = '  -1.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.00'
$ws.Range("E25").Value = '  -0.58%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("E27").Value = '  -1.80%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.111'
$ws.Range("E28").Value = '  -1.39%  '

$ws.Range("E29").Value = '  -1.53%  '

$ws.Range("E30").Value = '  -2.77%  '

$ws.Range("E31").Value = '  +0.05%  '

$ws.Range("E32").Value = '  -2.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.411.08'
$ws.Range("E33").Value = '  +8.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.93'
$ws.Range("E34").Value = '  -1.91%  '

$ws.Range("E36").Value = '  -1.34%  '

$ws.Range("E37").Value = '  -4.43%  '

$ws.Range("E39").Value = '  +1.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.86'
$ws.Range("E40").Value = '  +4.54%  '

$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.947'
$ws.Range("E42").Value = '  -13.21%  '

$ws.Range("E43").Value = '  +0.29%  '

$ws.Range("E44").Value = '  -0.60%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.723.74'
$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.80'
$ws.Range("E46").Value = '  -2.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.56'
$ws.Range("E47").Value = '  -2.48%  '

$ws.Range("E48").Value = '  -1.78%  '

$ws.Range("E49").Value = '  -0.45%  '

$ws.Range("E50").Value = '  -0.85%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0951'
$ws.Range("E51").Value = '  -3.29%  '
